$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes (col C, D, H) ---
# ColumnWidth as read back from the saved file is (input + 5/6), so we
# subtract 5/6 here to land on exact integer widths of 33 / 70 / 55.
$ws.Columns.Item(3).ColumnWidth = 32.16666666666667
$ws.Columns.Item(4).ColumnWidth = 69.16666666666667
$ws.Columns.Item(8).ColumnWidth = 54.16666666666667

# --- Row data (rows 2-14) ---
# Column A holds opportunity IDs that look numeric; a leading apostrophe
# forces them to be stored as text (matching the source data's text type)
# instead of being auto-converted to numbers.

# Row 2
$ws.Range("A2").Formula = "'1326535"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1326535"
$ws.Range("C2").Value = "ACCOUNTANT"
$ws.Range("D2").Value = "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "0 applicants"
$ws.Range("G2").Value = "3 - 6 Months"
$ws.Range("H2").Value = "Egypt holiday travel"

# Row 3
$ws.Range("A3").Formula = "'1326503"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1326503"
$ws.Range("C3").Value = "Data Analyst"
$ws.Range("D3").Value = "Αθήνα, Ελλάδα"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "6 applicants"
$ws.Range("G3").Value = "9 - 12 Weeks"
$ws.Range("H3").Value = "Inzeb"

# Row 4
$ws.Range("A4").Formula = "'1326324"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1326324"
$ws.Range("C4").Value = "International Relations Intern"
$ws.Range("D4").Value = "Jalandhar, Punjab, India"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "0 applicants"
$ws.Range("G4").Value = "3 - 6 Months"
$ws.Range("H4").Value = "Lovely Professional University"

# Row 5
$ws.Range("A5").Formula = "'1326081"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1326081"
$ws.Range("C5").Value = "Sales & Marketing"
$ws.Range("D5").Value = "Denizli, Kumkısık, Denizli, Türkiye"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "53 applicants"
$ws.Range("G5").Value = "6 - 18 Months"
$ws.Range("H5").Value = "ASM Crane"

# Row 6
$ws.Range("A6").Formula = "'1325826"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1325826"
$ws.Range("C6").Value = "Digital Marketing"
$ws.Range("D6").Value = "Cairo, Cairo Governorate, Egypt"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "0 applicants"
$ws.Range("G6").Value = "3 - 6 Months"
$ws.Range("H6").Value = "Kaian Arabi"

# Row 7
$ws.Range("A7").Formula = "'1317258"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1317258"
$ws.Range("C7").Value = "Sales Manager"
$ws.Range("D7").Value = "Bursa, Türkiye"
$ws.Range("E7").Value = "No"
$ws.Range("F7").Value = "104 applicants"
$ws.Range("G7").Value = "6 - 18 Months"
$ws.Range("H7").Value = "MAKELPORT MAKİNA VE YEDEK PARÇA"

# Row 8
$ws.Range("A8").Formula = "'1317128"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1317128"
$ws.Range("C8").Value = "BUSINESS ADMINISTRATION"
$ws.Range("D8").Value = "İstanbul, Türkiye"
$ws.Range("E8").Value = "No"
$ws.Range("F8").Value = "135 applicants"
$ws.Range("G8").Value = "6 - 18 Months"
$ws.Range("H8").Value = "Abdi İbrahim Pharmaceutical"

# Row 9
$ws.Range("A9").Formula = "'1312624"
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1312624"
$ws.Range("C9").Value = "MARKET RESEARCH/ANALYST"
$ws.Range("D9").Value = "Bursa, Türkiye"
$ws.Range("E9").Value = "No"
$ws.Range("F9").Value = "80 applicants"
$ws.Range("G9").Value = "6 - 18 Months"
$ws.Range("H9").Value = "YAŞALAR KALIP YEDEK PARÇA SANAYİ VE TİCARET LTD.ŞTİ."

# Row 10
$ws.Range("A10").Formula = "'1305678"
$ws.Range("B10").Value = "https://aiesec.org/opportunity/global-talent/1305678"
$ws.Range("C10").Value = "Market Research Analyst"
$ws.Range("D10").Value = "Bursa, Türkiye"
$ws.Range("E10").Value = "No"
$ws.Range("F10").Value = "121 applicants"
$ws.Range("G10").Value = "6 - 18 Months"
$ws.Range("H10").Value = "SİMYA GRUP MAKİNA"

# Row 11 (new)
$ws.Range("A11").Formula = "'1304097"
$ws.Range("B11").Value = "https://aiesec.org/opportunity/global-talent/1304097"
$ws.Range("C11").Value = "Marketing"
$ws.Range("D11").Value = "Bursa, Türkiye"
$ws.Range("E11").Value = "No"
$ws.Range("F11").Value = "79 applicants"
$ws.Range("G11").Value = "6 - 18 Months"
$ws.Range("H11").Value = "Tekinsan Otomotiv"

# Row 12 (new)
$ws.Range("A12").Formula = "'1301829"
$ws.Range("B12").Value = "https://aiesec.org/opportunity/global-talent/1301829"
$ws.Range("C12").Value = "Market Analyst"
$ws.Range("D12").Value = "Yıldırım, Türkiye"
$ws.Range("E12").Value = "No"
$ws.Range("F12").Value = "165 applicants"
$ws.Range("G12").Value = "6 - 18 Months"
$ws.Range("H12").Value = "Bemis Teknik Elektrik"

# Row 13 (new) - PREMIUM cell E13 is highlighted yellow
$ws.Range("A13").Formula = "'1299952"
$ws.Range("B13").Value = "https://aiesec.org/opportunity/global-talent/1299952"
$ws.Range("C13").Value = "Digital Marketing"
$ws.Range("D13").Value = "Yıldırım, Türkiye"
$ws.Range("E13").Value = "Yes"
$ws.Range("F13").Value = "88 applicants"
$ws.Range("G13").Value = "3 - 6 Months"
$ws.Range("H13").Value = "Via Premium"
$ws.Range("E13").Interior.Color = 65535

# Row 14 (new)
$ws.Range("A14").Formula = "'1294657"
$ws.Range("B14").Value = "https://aiesec.org/opportunity/global-talent/1294657"
$ws.Range("C14").Value = "Marketing"
$ws.Range("D14").Value = "Yıldırım, Türkiye"
$ws.Range("E14").Value = "No"
$ws.Range("F14").Value = "236 applicants"
$ws.Range("G14").Value = "6 - 18 Months"
$ws.Range("H14").Value = "OMSA"
